# Weight Loss Journey tracker - append the latest day's entry
# (date 2024-02-04 / serial 45326, weight 214) and scroll the sheet
# down to show the new row, matching how Excel left the view after
# the user typed in the new data at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right after the current last row (81) of the table.
$lastRow = 81
$newRow = $lastRow + 1

# Copy the date cell's formatting (short-date number format) down
# onto the new row so the new date cell matches the existing column
# instead of picking up Excel's default General format.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A" + $newRow).Value = 45326
$ws.Range("B" + $newRow).Value = 214

# Leave the selection where the user would naturally land next
# (the empty cell right below the newly-entered weight).
[void]$ws.Range("B" + ($newRow + 1)).Select()
